# Commit: Update countries & provincias Spain
# The "Pais" COVID dashboard data refreshed: Consejo Danes para los Refugiados and
# El Salvador received updated totals which bumped their rank in the (descending by
# total cases) table, shifting the rows in between down by one; Alemania, Ucrania and
# Uzbekistan got fresh totals without changing rank; Seychelles and Montserrat (tied on
# totals) swapped order; and the "last updated" timestamp banner advanced.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# "Datos actualizados a ..." timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 6 de Mayo de 2020 a las 09:03"

# Alemania: refreshed active/recovered counts
$ws.Cells.Item(9, 1).Value = "Alemania"
$ws.Cells.Item(9, 2).Value = 167007
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 4).Value = 137400
$ws.Cells.Item(9, 5).Value = 22614
$ws.Cells.Item(9, 6).Value = 1937
$ws.Cells.Item(9, 7).Value = 0
$ws.Cells.Item(9, 8).Value = 6993

# Ucrania: refreshed counts
$ws.Cells.Item(38, 1).Value = "Ucrania"
$ws.Cells.Item(38, 2).Value = 13184
$ws.Cells.Item(38, 3).Value = 487
$ws.Cells.Item(38, 4).Value = 2097
$ws.Cells.Item(38, 5).Value = 10760
$ws.Cells.Item(38, 6).Value = 169
$ws.Cells.Item(38, 7).Value = 11
$ws.Cells.Item(38, 8).Value = 327

# Uzbekistan: refreshed active/recovered counts
$ws.Cells.Item(72, 1).Value = "Uzbekistan"
$ws.Cells.Item(72, 2).Value = 2217
$ws.Cells.Item(72, 3).Value = 10
$ws.Cells.Item(72, 4).Value = 1539
$ws.Cells.Item(72, 5).Value = 668
$ws.Cells.Item(72, 6).Value = 8
$ws.Cells.Item(72, 7).Value = 0
$ws.Cells.Item(72, 8).Value = 10

# Consejo Danes para los Refugiados moves above Sudan (rows 99-115 shift down one rank)
$ws.Cells.Item(99, 1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(99, 2).Value = 797
$ws.Cells.Item(99, 3).Value = 92
$ws.Cells.Item(99, 4).Value = 92
$ws.Cells.Item(99, 5).Value = 670
$ws.Cells.Item(99, 6).Value = 0
$ws.Cells.Item(99, 7).Value = 1
$ws.Cells.Item(99, 8).Value = 35

$ws.Cells.Item(100, 1).Value = "Sudan"
$ws.Cells.Item(100, 2).Value = 778
$ws.Cells.Item(100, 3).Value = 0
$ws.Cells.Item(100, 4).Value = 70
$ws.Cells.Item(100, 5).Value = 663
$ws.Cells.Item(100, 6).Value = 0
$ws.Cells.Item(100, 7).Value = 0
$ws.Cells.Item(100, 8).Value = 45

$ws.Cells.Item(101, 1).Value = "Sri Lanka"
$ws.Cells.Item(101, 2).Value = 771
$ws.Cells.Item(101, 3).Value = 0
$ws.Cells.Item(101, 4).Value = 213
$ws.Cells.Item(101, 5).Value = 549
$ws.Cells.Item(101, 6).Value = 1
$ws.Cells.Item(101, 7).Value = 0
$ws.Cells.Item(101, 8).Value = 9

$ws.Cells.Item(102, 1).Value = "Guatemala"
$ws.Cells.Item(102, 2).Value = 763
$ws.Cells.Item(102, 3).Value = 33
$ws.Cells.Item(102, 4).Value = 79
$ws.Cells.Item(102, 5).Value = 665
$ws.Cells.Item(102, 6).Value = 5
$ws.Cells.Item(102, 7).Value = 0
$ws.Cells.Item(102, 8).Value = 19

$ws.Cells.Item(103, 1).Value = "Niger"
$ws.Cells.Item(103, 2).Value = 763
$ws.Cells.Item(103, 3).Value = 0
$ws.Cells.Item(103, 4).Value = 543
$ws.Cells.Item(103, 5).Value = 182
$ws.Cells.Item(103, 6).Value = 0
$ws.Cells.Item(103, 7).Value = 0
$ws.Cells.Item(103, 8).Value = 38

$ws.Cells.Item(104, 1).Value = "Costa Rica"
$ws.Cells.Item(104, 2).Value = 755
$ws.Cells.Item(104, 3).Value = 0
$ws.Cells.Item(104, 4).Value = 413
$ws.Cells.Item(104, 5).Value = 336
$ws.Cells.Item(104, 6).Value = 5
$ws.Cells.Item(104, 7).Value = 0
$ws.Cells.Item(104, 8).Value = 6

$ws.Cells.Item(105, 1).Value = "Principado de Andorra"
$ws.Cells.Item(105, 2).Value = 751
$ws.Cells.Item(105, 3).Value = 0
$ws.Cells.Item(105, 4).Value = 514
$ws.Cells.Item(105, 5).Value = 191
$ws.Cells.Item(105, 6).Value = 16
$ws.Cells.Item(105, 7).Value = 0
$ws.Cells.Item(105, 8).Value = 46

$ws.Cells.Item(106, 1).Value = "Libano"
$ws.Cells.Item(106, 2).Value = 741
$ws.Cells.Item(106, 3).Value = 0
$ws.Cells.Item(106, 4).Value = 206
$ws.Cells.Item(106, 5).Value = 510
$ws.Cells.Item(106, 6).Value = 43
$ws.Cells.Item(106, 7).Value = 0
$ws.Cells.Item(106, 8).Value = 25

$ws.Cells.Item(107, 1).Value = "Mayotte"
$ws.Cells.Item(107, 2).Value = 739
$ws.Cells.Item(107, 3).Value = 0
$ws.Cells.Item(107, 4).Value = 352
$ws.Cells.Item(107, 5).Value = 378
$ws.Cells.Item(107, 6).Value = 6
$ws.Cells.Item(107, 7).Value = 0
$ws.Cells.Item(107, 8).Value = 9

$ws.Cells.Item(108, 1).Value = "Crucero"
$ws.Cells.Item(108, 2).Value = 712
$ws.Cells.Item(108, 3).Value = 0
$ws.Cells.Item(108, 4).Value = 645
$ws.Cells.Item(108, 5).Value = 54
$ws.Cells.Item(108, 6).Value = 4
$ws.Cells.Item(108, 7).Value = 0
$ws.Cells.Item(108, 8).Value = 13

$ws.Cells.Item(109, 1).Value = "Burkina Faso"
$ws.Cells.Item(109, 2).Value = 688
$ws.Cells.Item(109, 3).Value = 0
$ws.Cells.Item(109, 4).Value = 548
$ws.Cells.Item(109, 5).Value = 92
$ws.Cells.Item(109, 6).Value = 0
$ws.Cells.Item(109, 7).Value = 0
$ws.Cells.Item(109, 8).Value = 48

$ws.Cells.Item(110, 1).Value = "Uruguay"
$ws.Cells.Item(110, 2).Value = 670
$ws.Cells.Item(110, 3).Value = 0
$ws.Cells.Item(110, 4).Value = 468
$ws.Cells.Item(110, 5).Value = 185
$ws.Cells.Item(110, 6).Value = 10
$ws.Cells.Item(110, 7).Value = 0
$ws.Cells.Item(110, 8).Value = 17

$ws.Cells.Item(111, 1).Value = "El Salvador"
$ws.Cells.Item(111, 2).Value = 633
$ws.Cells.Item(111, 3).Value = 46
$ws.Cells.Item(111, 4).Value = 219
$ws.Cells.Item(111, 5).Value = 400
$ws.Cells.Item(111, 6).Value = 4
$ws.Cells.Item(111, 7).Value = 1
$ws.Cells.Item(111, 8).Value = 14

$ws.Cells.Item(112, 1).Value = "Mali"
$ws.Cells.Item(112, 2).Value = 612
$ws.Cells.Item(112, 3).Value = 0
$ws.Cells.Item(112, 4).Value = 228
$ws.Cells.Item(112, 5).Value = 352
$ws.Cells.Item(112, 6).Value = 0
$ws.Cells.Item(112, 7).Value = 0
$ws.Cells.Item(112, 8).Value = 32

$ws.Cells.Item(113, 1).Value = "Georgia"
$ws.Cells.Item(113, 2).Value = 610
$ws.Cells.Item(113, 3).Value = 6
$ws.Cells.Item(113, 4).Value = 269
$ws.Cells.Item(113, 5).Value = 332
$ws.Cells.Item(113, 6).Value = 6
$ws.Cells.Item(113, 7).Value = 0
$ws.Cells.Item(113, 8).Value = 9

$ws.Cells.Item(114, 1).Value = "San Marino"
$ws.Cells.Item(114, 2).Value = 589
$ws.Cells.Item(114, 3).Value = 0
$ws.Cells.Item(114, 4).Value = 92
$ws.Cells.Item(114, 5).Value = 456
$ws.Cells.Item(114, 6).Value = 5
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 41

$ws.Cells.Item(115, 1).Value = "Maldivas"
$ws.Cells.Item(115, 2).Value = 573
$ws.Cells.Item(115, 3).Value = 0
$ws.Cells.Item(115, 4).Value = 20
$ws.Cells.Item(115, 5).Value = 551
$ws.Cells.Item(115, 6).Value = 2
$ws.Cells.Item(115, 7).Value = 0
$ws.Cells.Item(115, 8).Value = 2

# Seychelles / Montserrat swap order (tie on total cases)
$ws.Cells.Item(205, 1).Value = "Montserrat"
$ws.Cells.Item(205, 2).Value = 11
$ws.Cells.Item(205, 3).Value = 0
$ws.Cells.Item(205, 4).Value = 7
$ws.Cells.Item(205, 5).Value = 3
$ws.Cells.Item(205, 6).Value = 1
$ws.Cells.Item(205, 7).Value = 0
$ws.Cells.Item(205, 8).Value = 1

$ws.Cells.Item(206, 1).Value = "Seychelles"
$ws.Cells.Item(206, 2).Value = 11
$ws.Cells.Item(206, 3).Value = 0
$ws.Cells.Item(206, 4).Value = 8
$ws.Cells.Item(206, 5).Value = 3
$ws.Cells.Item(206, 6).Value = 0
$ws.Cells.Item(206, 7).Value = 0
$ws.Cells.Item(206, 8).Value = 0

